# Refresh the cryptocurrency price/volume snapshot (and fix two rank swaps)
# as published by the upstream GitHub Actions scraper job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text like '70.946.15' or '0.999' -- these look numeric
# to Excel's auto-detection, so force Text formatting before writing, then
# drop the temporary format again so the cell style matches the original file.
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

# Row 2
Set-TextValue "D2" '71.038.99'
$ws.Range("E2").Value = '  -1.91%  '

# Row 3
Set-TextValue "D3" '3.947.26'
$ws.Range("E3").Value = '  -2.47%  '

# Row 4
Set-TextValue "D4" '1.00'
$ws.Range("E4").Value = '  -0.10%  '

# Row 5
Set-TextValue "D5" '536.32'
$ws.Range("E5").Value = '  +2.73%  '

# Row 6
Set-TextValue "D6" '147.91'
$ws.Range("E6").Value = '  -0.48%  '

# Row 7
Set-TextValue "D7" '3.942.82'
$ws.Range("E7").Value = '  -2.39%  '

# Row 8
Set-TextValue "D8" '0.688'
$ws.Range("E8").Value = '  -4.47%  '

# Row 9
Set-TextValue "D9" '1.00'
$ws.Range("E9").Value = '  -0.02%  '

# Row 10
Set-TextValue "D10" '0.738'
$ws.Range("E10").Value = '  -5.55%  '

# Row 11
$ws.Range("E11").Value = '  -7.63%  '

# Row 12
Set-TextValue "D12" '55.17'
$ws.Range("E12").Value = '  +13.91%  '

# Row 13
Set-TextValue "D13" '0.0000317'
$ws.Range("E13").Value = '  -5.15%  '

# Row 14
Set-TextValue "D14" '10.60'
$ws.Range("E14").Value = '  -5.28%  '

# Row 15
Set-TextValue "D15" '4.570.10'
$ws.Range("E15").Value = '  -2.65%  '

# Row 16
Set-TextValue "D16" '3.943.01'
$ws.Range("E16").Value = '  -2.53%  '

# Row 17
$ws.Range("B17").Value = 'Uniswap'
$ws.Range("C17").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue "D17" '13.88'
$ws.Range("E17").Value = '  -3.21%  '

# Row 18
$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue "D18" '20.49'
$ws.Range("E18").Value = '  -3.92%  '

# Row 20
Set-TextValue "D20" '1.17'
$ws.Range("E20").Value = '  -4.85%  '

# Row 21
Set-TextValue "D21" '70.896.13'
$ws.Range("E21").Value = '  -2.01%  '

# Row 22
Set-TextValue "D22" '421.25'
$ws.Range("E22").Value = '  -5.71%  '

# Row 23
Set-TextValue "D23" '3.61'
$ws.Range("E23").Value = '  -0.11%  '

# Row 24
Set-TextValue "D24" '97.32'
$ws.Range("E24").Value = '  -7.35%  '

# Row 25
Set-TextValue "D25" '4.22'
$ws.Range("E25").Value = '  +4.03%  '

# Row 26
Set-TextValue "D26" '14.45'
$ws.Range("E26").Value = '  -4.99%  '

# Row 27
Set-TextValue "D27" '11.40'
$ws.Range("E27").Value = '  -1.25%  '

# Row 28
Set-TextValue "D28" '3.82'
$ws.Range("E28").Value = '  +16.10%  '

# Row 29
Set-TextValue "D29" '10.68'
$ws.Range("E29").Value = '  -4.63%  '

# Row 30
Set-TextValue "D30" '5.87'
$ws.Range("E30").Value = '  +0.84%  '

# Row 31
Set-TextValue "D31" '36.43'
$ws.Range("E31").Value = '  -4.51%  '

# Row 32
Set-TextValue "D32" '7.91'
$ws.Range("E32").Value = '  +17.45%  '

# Row 33
Set-TextValue "D33" '50.79'
$ws.Range("E33").Value = '  +18.68%  '

# Row 34
$ws.Range("E34").Value = '  -0.14%  '

# Row 35
Set-TextValue "D35" '13.35'
$ws.Range("E35").Value = '  -3.65%  '

# Row 36
Set-TextValue "D36" '683.99'
$ws.Range("E36").Value = '  +0.19%  '

# Row 37
Set-TextValue "D37" '65.38'
$ws.Range("E37").Value = '  -3.79%  '

# Row 38
Set-TextValue "D38" '0.442'
$ws.Range("E38").Value = '  +2.17%  '

# Row 39
Set-TextValue "D39" '0.0₃0816'
$ws.Range("E39").Value = '  -7.17%  '

# Row 40
$ws.Range("E40").Value = '  -3.19%  '

# Row 41
Set-TextValue "D41" '3.37'
$ws.Range("E41").Value = '  -4.16%  '

# Row 42
Set-TextValue "D42" '0.999'
$ws.Range("E42").Value = '  -0.07%  '

# Row 43
Set-TextValue "D43" '0.999'
$ws.Range("E43").Value = '  +0.05%  '

# Row 44
Set-TextValue "D44" '0.0481'
$ws.Range("E44").Value = '  -4.17%  '

# Row 45
$ws.Range("E45").Value = '  -0.91%  '

# Row 46
$ws.Range("B46").Value = 'THORChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue "D46" '10.00'
$ws.Range("E46").Value = '  +2.31%  '

# Row 47
$ws.Range("B47").Value = 'Stellar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue "D47" '0.149'
$ws.Range("E47").Value = '  -5.90%  '

# Row 48
Set-TextValue "D48" '2.67'

# Row 49
$ws.Range("E49").Value = '  -2.91%  '

# Row 50
$ws.Range("E50").Value = '  -2.61%  '

# Row 51
Set-TextValue "D51" '144.63'
$ws.Range("E51").Value = '  -0.12%  '
